# Updates cryptos list figures (price/volume columns) to the latest
# scraped snapshot, matching commit "Updated cryptos list ... with GitHub Actions".
# Numeric-looking text values in column D are prefixed with a leading
# apostrophe so Excel stores them as text (preserving formatting such as
# trailing zeros) instead of silently coercing them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.821.38'
$ws.Range("E2").Value = '  -1.27%  '
$ws.Range("D3").Value = '2.909.54'
$ws.Range("E3").Value = '  -1.75%  '
$ws.Range("D5").Value = '''586.73'
$ws.Range("E5").Value = '  -1.55%  '
$ws.Range("D6").Value = '''146.61'
$ws.Range("E6").Value = '  +0.97%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +0.27%  '
$ws.Range("D9").Value = '2.908.43'
$ws.Range("D10").Value = '''6.88'
$ws.Range("E10").Value = '  -6.41%  '
$ws.Range("E11").Value = '  +4.80%  '
$ws.Range("D12").Value = '''0.435'
$ws.Range("E12").Value = '  -2.61%  '
$ws.Range("E13").Value = '  +1.62%  '
$ws.Range("D14").Value = '''32.85'
$ws.Range("E14").Value = '  -1.56%  '
$ws.Range("E15").Value = '  -0.80%  '
$ws.Range("D16").Value = '3.391.30'
$ws.Range("E16").Value = '  -1.62%  '
$ws.Range("D17").Value = '61.854.17'
$ws.Range("E17").Value = '  -1.00%  '
$ws.Range("D18").Value = '''6.64'
$ws.Range("E18").Value = '  -1.12%  '
$ws.Range("D19").Value = '2.905.83'
$ws.Range("E19").Value = '  -1.73%  '
$ws.Range("D20").Value = '''436.98'
$ws.Range("E20").Value = '  -1.00%  '
$ws.Range("D21").Value = '''13.39'
$ws.Range("E21").Value = '  -0.30%  '
$ws.Range("D22").Value = '''0.661'
$ws.Range("E22").Value = '  -1.63%  '
$ws.Range("E23").Value = '  -2.06%  '
$ws.Range("D24").Value = '''81.11'
$ws.Range("E24").Value = '  -0.88%  '
$ws.Range("D25").Value = '''11.98'
$ws.Range("E25").Value = '  -0.10%  '
$ws.Range("E26").Value = '  -7.55%  '
$ws.Range("D27").Value = '''2.08'
$ws.Range("E27").Value = '  -2.38%  '
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("D29").Value = '''0.0000106'
$ws.Range("E29").Value = '  +21.18%  '
$ws.Range("D30").Value = '''7.16'
$ws.Range("E30").Value = '  +1.48%  '
$ws.Range("E31").Value = '  -1.72%  '
$ws.Range("E32").Value = '  -0.64%  '
$ws.Range("E33").Value = '  +0.46%  '
$ws.Range("E34").Value = '  -0.06%  '
$ws.Range("D35").Value = '''25.91'
$ws.Range("E35").Value = '  -2.58%  '
$ws.Range("D36").Value = '''0.972'
$ws.Range("E36").Value = '  -1.98%  '
$ws.Range("D37").Value = '''5.51'
$ws.Range("E37").Value = '  -2.17%  '
$ws.Range("D38").Value = '''3.03'
$ws.Range("E38").Value = '  +3.82%  '
$ws.Range("D39").Value = '''49.20'
$ws.Range("E39").Value = '  -0.82%  '
$ws.Range("D40").Value = '''2.00'
$ws.Range("E40").Value = '  -2.92%  '
$ws.Range("D41").Value = '''8.39'
$ws.Range("E41").Value = '  -1.97%  '
$ws.Range("E42").Value = '  -1.07%  '
$ws.Range("E43").Value = '  -3.27%  '
$ws.Range("D44").Value = '''38.99'
$ws.Range("E44").Value = '  -0.27%  '
$ws.Range("D45").Value = '2.701.95'
$ws.Range("E45").Value = '  -0.57%  '
$ws.Range("D46").Value = '''133.90'
$ws.Range("E46").Value = '  -1.16%  '
$ws.Range("E47").Value = '  -1.20%  '
$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D48").Value = '''342.52'
$ws.Range("E48").Value = '  -5.48%  '
$ws.Range("B49").Value = 'USDe'
$ws.Range("C49").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D49").Value = '''1.00'
$ws.Range("E49").Value = '  +0.01%  '
$ws.Range("E50").Value = '  -1.51%  '
$ws.Range("D51").Value = '''22.32'
$ws.Range("E51").Value = '  -2.60%  '
